$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update stimulus filename cells (row 2, columns C-F) to new restructured paths
$ws.Range("C2").Value = "stimuli/Social1.png"
$ws.Range("D2").Value = "stimuli/Social2.png"
$ws.Range("E2").Value = "stimuli/Nonsocial1.png"
$ws.Range("F2").Value = "stimuli/Nonsocial2.png"

# Set column widths to reflect the new (shorter) content widths
$ws.Range("C1:D1").ColumnWidth = 15.666666666666666
$ws.Range("E1:F1").ColumnWidth = 19.0

# Update the active selection to G2
$ws.Range("G2").Select()
